$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 24 (pushes the blank rows + signature block down by one)
$ws.Rows("24").Insert()

# 2) Copy the "closing" (bottom-border) row format that used to belong to row 23
#    down onto the brand-new row 24, before we overwrite row 23's own format.
$ws.Range("B23:J23").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)

# 3) Row 23 is no longer the last data row, so give it the regular "middle" row
#    formatting (same as the row above it).
$ws.Range("B22:J22").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4) Populate the new row 24 with the additional period (2509) for the same worker.
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "30763004"
$ws.Range("D24").Value = "ANA ELVIRA HERNANDEZ DE BARRIOS"
$ws.Range("E24").Value = "2509"
$ws.Range("F24").Value = 52000
$ws.Range("G24").Value = 1300000

# 5) The "Periodo Mora" column is centered for every data row, including the new one.
$ws.Range("E16:E24").HorizontalAlignment = -4108

# 6) Update the summary figures: total overdue value and number of periods.
$ws.Range("E11").Value = 468000
$ws.Range("F13").Value = 9
